$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Info sheet: the "Exam Type Code" / "MIDTERM_1/2025" row is removed, the
# "Exam" label becomes "Exam Type", and every row below shifts up by one
# (columns A:C only -- the P:R lookup table used by the VLOOKUP formulas
# stays put).
# -----------------------------------------------------------------------
$info = $wb.Worksheets.Item("Info")

# "Exam" -> "Exam Type" (row 3 keeps its VLOOKUP formula in B3)
$info.Range("A3").Value = "Exam Type"

# Drop the old "Exam Type Code" / "MIDTERM_1/2025" row entirely.
$info.Range("A4:C4").Clear()

# Shift each remaining label/value pair up one row, left-to-right,
# top-to-bottom so a source row is always read before it is overwritten.
$info.Range("A6:B6").Cut($info.Range("A5"))
$info.Range("A6:B6").Clear()

$info.Range("A7:B7").Cut($info.Range("A6"))
$info.Range("A7:B7").Clear()

$info.Range("A9:B9").Cut($info.Range("A8"))
$info.Range("A9:B9").Clear()

$info.Range("A10:B10").Cut($info.Range("A9"))
$info.Range("A10:B10").Clear()

$info.Range("A11:B11").Cut($info.Range("A10"))
$info.Range("A11:B11").Clear()

$info.Range("A13:B13").Cut($info.Range("A12"))
$info.Range("A13:B13").Clear()

$info.Range("A14:C14").Cut($info.Range("A13"))
$info.Range("A14:C14").Clear()

# Cut() only carries the cached value, not the formula -- restore the
# formulas on the cells that now live one row higher.
$info.Range("B12").Formula = "=DATE(2025,6,10)"
$info.Range("B13").Formula = "=TIME(9,0,0)"
$info.Range("C13").Formula = "=TIME(12,0,0)"

# The selected cell moved too.
[void]$info.Range("F11").Select()

# -----------------------------------------------------------------------
# WrittenQuestion sheet: header/question cells lose/gain the "top
# aligned" style that the first two rows used to use exclusively.
# -----------------------------------------------------------------------
$written = $wb.Worksheets.Item("WrittenQuestion")

# A2/A3 no longer use the wrap-text "top aligned" style -- back to Normal.
$written.Range("A2").Style = "Normal"
$written.Range("A3").Style = "Normal"

# B4/D4 pick up the "top aligned" (no wrap) style already used by C4.
$written.Range("B4").VerticalAlignment = -4160
$written.Range("B4").WrapText = $false
$written.Range("D4").VerticalAlignment = -4160
$written.Range("D4").WrapText = $false
